$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Mobile column (F) is formatted as Text so numeric-looking
# values are stored as strings, matching the source data.
$ws.Range("F3:F4").NumberFormat = "@"

# Update row 3 (existing candidate record)
$ws.Range("B3").Value = 312
$ws.Range("C3").Value = "aman"
$ws.Range("D3").Value = "aman@gamil.com"
$ws.Range("E3").Value = "tcs"
$ws.Range("F3").Value = "564321789"
$ws.Range("G3").Value = "DM_selected"

# Add new row 4 (new candidate record)
$ws.Range("B4").Value = 313
$ws.Range("C4").Value = "abcdf"
$ws.Range("D4").Value = "tishya@gmail.com"
$ws.Range("E4").Value = "globalTiger"
$ws.Range("F4").Value = "55667788"
$ws.Range("G4").Value = "DM_selected"
